# Refresh the "cryptos" price/volume table with the latest scraped values.
# Columns: D = Price (text), E = Volume(1h) change (text, padded with spaces).
# Some new Price values look like plain decimals (e.g. "578.98"), so for those
# cells we briefly force a Text number format before assigning the value -
# otherwise Excel's input parser would silently convert them to numbers -
# then restore the default "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.454.18'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '3.235.01'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '183.86'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.33%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.608'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.12%  '
$ws.Range('D9').Value = '3.233.65'
$ws.Range('E9').Value = '  -1.25%  '
$ws.Range('E10').Value = '  -3.47%  '
$ws.Range('E11').Value = '  -2.38%  '
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').Value = '3.793.90'
$ws.Range('E13').Value = '  -1.29%  '
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.64'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.61%  '
$ws.Range('D16').Value = '67.533.42'
$ws.Range('E16').Value = '  -0.96%  '
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('D18').Value = '3.223.27'
$ws.Range('E18').Value = '  -1.56%  '
$ws.Range('E19').Value = '  -1.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '394.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.55'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.43%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  -2.77%  '
$ws.Range('E27').Value = '  -1.55%  '
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.96'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.10%  '
$ws.Range('E31').Value = '  -4.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '22.61'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E34').Value = '  -2.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.79'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.41%  '
$ws.Range('E37').Value = '  -4.06%  '
$ws.Range('E38').Value = '  +0.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.46'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('E40').Value = '  -4.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.54'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.65%  '
$ws.Range('E42').Value = '  -5.07%  '
$ws.Range('E43').Value = '  -7.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0687'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('D45').Value = '2.613.52'
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.50'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '24.72'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '333.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0278'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.31'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.77%  '
$ws.Range('E51').Value = '  -0.78%  '
